$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (avoid Excel auto-converting numeric-looking strings to numbers and
# dropping formatting such as trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.289.63"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.333.73"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "548.42"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").Value = "131.22"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.580"
$ws.Range("E8").Value = "  -1.32%  "
$ws.Range("D9").Value = "2.333.20"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "0.337"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "23.72"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "2.754.64"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "60.206.89"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "2.323.90"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "4.10"
$ws.Range("E20").Value = "  -1.56%  "
$ws.Range("D21").Value = "314.73"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "6.59"
$ws.Range("E22").Value = "  -3.57%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "64.16"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "7.95"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "1.26"
$ws.Range("E29").Value = "  +7.62%  "
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "171.27"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").Value = "6.11"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D34").Value = "0.385"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.37"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "4.10"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "323.97"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "38.15"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "137.77"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("D44").Value = "3.51"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "19.34"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Value = "0.566"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "0.0499"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "0.0₆0219"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("D51").Value = "10.94"
$ws.Range("E51").Value = "  -0.75%  "

Write-Output "Applied 90 cell updates"
